# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the QoQ forecast-error summary table (rows Q1..Q9) with
# recalculated ME / MAE / MSE / RMSE / SE / N values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.0725292808699395
$ws.Range("C2").Value = 1.452525492729129
$ws.Range("D2").Value = 3.943704826413877
$ws.Range("E2").Value = 1.985876337140326
$ws.Range("F2").Value = 2.004298684434776
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.970128790070132
$ws.Range("C3").Value = 1.91197192143461
$ws.Range("D3").Value = 11.00797031802001
$ws.Range("E3").Value = 3.317826143428858
$ws.Range("F3").Value = 3.205037560462687
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 1.027439519176732
$ws.Range("C4").Value = 2.670308812046868
$ws.Range("D4").Value = 20.8551199302896
$ws.Range("E4").Value = 4.566740624372004
$ws.Range("F4").Value = 4.495773270230817
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 1.040770019679559
$ws.Range("C5").Value = 2.551126684817124
$ws.Range("D5").Value = 21.80091093632453
$ws.Range("E5").Value = 4.669144561514939
$ws.Range("F5").Value = 4.599838150987386
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.6404470419866763
$ws.Range("C6").Value = 2.649364789166133
$ws.Range("D6").Value = 22.40911179286584
$ws.Range("E6").Value = 4.733826337421541
$ws.Range("F6").Value = 4.741010144670577
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.6619285698089534
$ws.Range("C7").Value = 2.518113852960962
$ws.Range("D7").Value = 23.16329070585441
$ws.Range("E7").Value = 4.812825646733363
$ws.Range("F7").Value = 4.83108003263652
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.6925657535171938
$ws.Range("C8").Value = 2.61980777971489
$ws.Range("D8").Value = 24.31016614921836
$ws.Range("E8").Value = 4.930534063285474
$ws.Range("F8").Value = 4.948987541160202
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = 0.7380237456049207
$ws.Range("C9").Value = 3.434955442878508
$ws.Range("D9").Value = 40.91359093038702
$ws.Range("E9").Value = 6.396373263841551
$ws.Range("F9").Value = 6.518710873632213
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -1.030652359722206
$ws.Range("C10").Value = 3.294535520879184
$ws.Range("D10").Value = 28.3461002497294
$ws.Range("E10").Value = 5.324105582135783
$ws.Range("F10").Value = 5.436681950116494
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = -0.1548825442147432
$ws.Range("C11").Value = 2.972445703758152
$ws.Range("D11").Value = 9.566860607955785
$ws.Range("E11").Value = 3.093034207369163
$ws.Range("F11").Value = 3.45377909062185
$ws.Range("G11").Value = 4

